# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update the "last updated" timestamp in the title cell (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 18 de Agosto de 2020 a las 16:55"

# --- Swap country names between row 213 and row 214 ---
# (Islas Malvinas / Montserrat exchange places in the shared-string list,
#  while their underlying stats stay tied to the row below)
$ws.Range("A213").Value = "Montserrat"
$ws.Range("A214").Value = "Islas Malvinas"

# --- Update per-country statistics (B:H = Casos totales, Nuevos casos,
#     Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes) ---

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 5614218
$ws.Range("C4").Value = 2191
$ws.Range("E4").Value = 2465641
$ws.Range("G4").Value = 73
$ws.Range("H4").Value = 173789

# Row 6 - Rusia
$ws.Range("B6").Value = 2732218
$ws.Range("C6").Value = 30614
$ws.Range("D6").Value = 2005215
$ws.Range("E6").Value = 674723
$ws.Range("G6").Value = 355
$ws.Range("H6").Value = 52280

# Row 22
$ws.Range("B22").Value = 227217
$ws.Range("C22").Value = 531
$ws.Range("E22").Value = 15017
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = 9300

# Row 50
$ws.Range("B50").Value = 54448
$ws.Range("C50").Value = 214
$ws.Range("D50").Value = 39936
$ws.Range("E50").Value = 12728
$ws.Range("G50").Value = 5
$ws.Range("H50").Value = 1784

# Row 59
$ws.Range("D59").Value = 33500
$ws.Range("E59").Value = 2957

# Row 65
$ws.Range("B65").Value = 30789
$ws.Range("C65").Value = 412
$ws.Range("E65").Value = 8661

# Row 66
$ws.Range("B66").Value = 30636
$ws.Range("C66").Value = 271
$ws.Range("D66").Value = 17368
$ws.Range("E66").Value = 12781
$ws.Range("G66").Value = 5
$ws.Range("H66").Value = 487

# Row 87
$ws.Range("B87").Value = 10089
$ws.Range("C87").Value = 29
$ws.Range("E87").Value = 970

# Row 112
$ws.Range("B112").Value = 4464
$ws.Range("C112").Value = 120
$ws.Range("D112").Value = 2407
$ws.Range("E112").Value = 2020
$ws.Range("G112").Value = 1
$ws.Range("H112").Value = 37

# Row 166
$ws.Range("B166").Value = 600
$ws.Range("C166").Value = 12
$ws.Range("E166").Value = 448

# Row 213 (now Montserrat)
$ws.Range("D213").Value = 12
$ws.Range("H213").Value = 1

# Row 214 (now Islas Malvinas)
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0
